# Apply updated odds values to Sheet1, as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 3.1
$ws.Range("L5").Value = 1.73
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 3.5
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 1.78
$ws.Range("Q5").Value = 2.03
$ws.Range("T5").Value = 5.5
$ws.Range("U5").Value = 10
$ws.Range("Z5").Value = 5
$ws.Range("AB5").Value = 23
$ws.Range("AC5").Value = 101
$ws.Range("AF5").Value = 13
$ws.Range("AG5").Value = 13

# Row 6
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 1.08
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 1.4
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = 2.15
$ws.Range("O6").Value = 1.67
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.63
$ws.Range("AG6").Value = 15

# Row 8
$ws.Range("G8").Value = 1.93
$ws.Range("I8").Value = 3.8
$ws.Range("L8").Value = 1.28
$ws.Range("M8").Value = 3.45
$ws.Range("N8").Value = 1.85
$ws.Range("O8").Value = 1.9
$ws.Range("P8").Value = 1.4
$ws.Range("Q8").Value = 2.8
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 2.02
$ws.Range("T8").Value = 7.6
$ws.Range("V8").Value = 8.75
$ws.Range("W8").Value = 18.5
$ws.Range("X8").Value = 16
$ws.Range("Y8").Value = 26
$ws.Range("AA8").Value = 7
$ws.Range("AC8").Value = 65
$ws.Range("AE8").Value = 10.75
$ws.Range("AF8").Value = 22
$ws.Range("AG8").Value = 13.5
$ws.Range("AH8").Value = 60
$ws.Range("AI8").Value = 37
$ws.Range("AJ8").Value = 45

# Row 9
$ws.Range("G9").Value = 2.22
$ws.Range("H9").Value = 3.15
$ws.Range("I9").Value = 3.3
$ws.Range("J9").Value = 1.1
$ws.Range("K9").Value = 6.3
$ws.Range("L9").Value = 1.44
$ws.Range("M9").Value = 2.65
$ws.Range("N9").Value = 2.3
$ws.Range("O9").Value = 1.57
$ws.Range("P9").Value = 1.5
$ws.Range("Q9").Value = 2.5
$ws.Range("T9").Value = 6.4
$ws.Range("U9").Value = 10.5
$ws.Range("W9").Value = 23
$ws.Range("Z9").Value = 6.3
$ws.Range("AA9").Value = 6.4
$ws.Range("AB9").Value = 18.5
$ws.Range("AE9").Value = 7.9
$ws.Range("AF9").Value = 17
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 50
$ws.Range("AI9").Value = 37

# Row 10
$ws.Range("L10").Value = 1.44
$ws.Range("M10").Value = 2.63

# Row 39
$ws.Range("J39").Value = 1.11
$ws.Range("K39").Value = 6.5
$ws.Range("N39").Value = 2.7
$ws.Range("O39").Value = 1.44

# Row 40
$ws.Range("G40").Value = 1.73
$ws.Range("H40").Value = 3.25
$ws.Range("I40").Value = 5.25
$ws.Range("X40").Value = 15
$ws.Range("AC40").Value = 67
$ws.Range("AF40").Value = 26
